$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VCRA")

# --- Insert two new quarterly columns before column D ---
# (existing D:K quarterly data shifts right to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy cell formatting (number format / font / alignment) from column F
# into the two freshly inserted columns D:E so the new quarters look like
# the rest of the table (dates keep the custom date format, numbers keep
# the #,##0 style).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the two new quarter columns (2018-12-31 = D, 2018-09-30 = E) ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 48900
$ws.Range("E8").Value = 47800
$ws.Range("D9").Value = 17700
$ws.Range("E9").Value = 16700
$ws.Range("D10").Value = 31200
$ws.Range("E10").Value = 31100
$ws.Range("D12").Value = 8200
$ws.Range("E12").Value = 8000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 48300
$ws.Range("E17").Value = 46800
$ws.Range("D18").Value = 600
$ws.Range("E18").Value = 1000
$ws.Range("D20").Value = 600
$ws.Range("E20").Value = 1000
$ws.Range("D21").Value = 3200
$ws.Range("E21").Value = 3900
$ws.Range("D22").Value = 2100
$ws.Range("E22").Value = 2100
$ws.Range("D23").Value = -1000
$ws.Range("E23").Value = -100
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -1100
$ws.Range("E26").Value = -200
$ws.Range("D27").Value = -1100
$ws.Range("E27").Value = -200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -600
$ws.Range("E32").Value = -1000
$ws.Range("D33").Value = -1100
$ws.Range("E33").Value = -200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -1100
$ws.Range("E35").Value = -200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 34300
$ws.Range("E41").Value = 33000
$ws.Range("D42").Value = 186900
$ws.Range("E42").Value = 182000
$ws.Range("D43").Value = 44300
$ws.Range("E43").Value = 39300
$ws.Range("D44").Value = 4400
$ws.Range("E44").Value = 3700
$ws.Range("D45").Value = 4700
$ws.Range("E45").Value = 4500
$ws.Range("D46").Value = 274500
$ws.Range("E46").Value = 262400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 7500
$ws.Range("E48").Value = 7100
$ws.Range("D49").Value = 58300
$ws.Range("E49").Value = 59500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 11800
$ws.Range("E52").Value = 11200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 352100
$ws.Range("E54").Value = 340100
$ws.Range("D57").Value = 4200
$ws.Range("E57").Value = 3200
$ws.Range("D58").Value = 1000
$ws.Range("E58").Value = 900
$ws.Range("D59").Value = 56000
$ws.Range("E59").Value = 54000
$ws.Range("D60").Value = 61200
$ws.Range("E60").Value = 58000
$ws.Range("D61").Value = 110500
$ws.Range("E61").Value = 108900
$ws.Range("D62").Value = 17500
$ws.Range("E62").Value = 15800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 189200
$ws.Range("E66").Value = 182800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -132300
$ws.Range("E72").Value = -131200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 162900
$ws.Range("E76").Value = 157400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -1100
$ws.Range("E81").Value = -200
$ws.Range("D83").Value = 2000
$ws.Range("E83").Value = 1900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 6400
$ws.Range("E89").Value = 6500
$ws.Range("D91").Value = -1300
$ws.Range("E91").Value = -1700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -6300
$ws.Range("E94").Value = -4800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 1200
$ws.Range("E100").Value = 3800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 1300
$ws.Range("E102").Value = 5500

# --- A handful of historical quarters were also corrected alongside the insert ---
$ws.Range("H8").Value = 44100
$ws.Range("I8").Value = 45600
$ws.Range("H10").Value = 28000
$ws.Range("I10").Value = 29200
$ws.Range("H17").Value = 44200
$ws.Range("I17").Value = 44000
$ws.Range("H18").Value = -100
$ws.Range("I18").Value = 1600
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 3700
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 1700
$ws.Range("H24").Value = -300
$ws.Range("H26").Value = 400
$ws.Range("I26").Value = 1400
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 1400
$ws.Range("H33").Value = 1100
$ws.Range("I33").Value = 1400
$ws.Range("H35").Value = 1100
$ws.Range("I35").Value = 1400
$ws.Range("H81").Value = 1100
$ws.Range("I81").Value = 1400
$ws.Range("J91").Value = -1100